$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/17/2025  Through  11/23/2025"

# --- Numeric cell updates ---
$ws.Range("N15").Value = 54.545454545454
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 70
$ws.Range("I16").Value = 144
$ws.Range("J16").Value = 146
$ws.Range("K16").Value = -1.369863013698
$ws.Range("L16").Value = -9.433962264150
$ws.Range("M16").Value = -33.944954128440
$ws.Range("N16").Value = 65.517241379310
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = 60
$ws.Range("I17").Value = 242
$ws.Range("J17").Value = 212
$ws.Range("K17").Value = 14.150943396226
$ws.Range("L17").Value = 6.140350877192
$ws.Range("M17").Value = 46.666666666666
$ws.Range("N17").Value = 290.322580645161
$ws.Range("C18").Value = 1
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -57.142857142857
$ws.Range("I18").Value = 94
$ws.Range("K18").Value = -10.476190476190
$ws.Range("L18").Value = -10.476190476190
$ws.Range("M18").Value = -3.092783505154
$ws.Range("N18").Value = -10.476190476190
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 3.571428571428
$ws.Range("I19").Value = 333
$ws.Range("J19").Value = 325
$ws.Range("K19").Value = 2.461538461538
$ws.Range("L19").Value = 2.777777777777
$ws.Range("M19").Value = 36.475409836065
$ws.Range("N19").Value = 516.666666666667
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 12.5
$ws.Range("I20").Value = 126
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 26
$ws.Range("L20").Value = -2.325581395348
$ws.Range("M20").Value = 103.225806451613
$ws.Range("N20").Value = 51.807228915662
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 76
$ws.Range("G21").Value = 65
$ws.Range("H21").Value = 16.923076923076
$ws.Range("I21").Value = 956
$ws.Range("J21").Value = 903
$ws.Range("K21").Value = 5.869324473975
$ws.Range("L21").Value = -0.312825860271
$ws.Range("M21").Value = 18.463444857496
$ws.Range("N21").Value = 137.810945273632
$ws.Range("F22").Value = 2
$ws.Range("I22").Value = 7
$ws.Range("K22").Value = -53.333333333333
$ws.Range("L22").Value = -74.074074074074
$ws.Range("M22").Value = -50
$ws.Range("C23").Value = 1
$ws.Range("I23").Value = 18
$ws.Range("K23").Value = -33.333333333333
$ws.Range("L23").Value = 28.571428571428
$ws.Range("M23").Value = 63.636363636363
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 8.695652173913
$ws.Range("F24").Value = 74
$ws.Range("G24").Value = 70
$ws.Range("H24").Value = 5.714285714285
$ws.Range("I24").Value = 730
$ws.Range("J24").Value = 831
$ws.Range("K24").Value = -12.154031287605
$ws.Range("L24").Value = -12.154031287605
$ws.Range("M24").Value = 36.960600375234
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -92.857142857142
$ws.Range("F25").Value = 7
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = -78.125
$ws.Range("I25").Value = 139
$ws.Range("J25").Value = 279
$ws.Range("K25").Value = -50.179211469534
$ws.Range("L25").Value = -27.225130890052
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -80
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = -16
$ws.Range("I26").Value = 329
$ws.Range("J26").Value = 309
$ws.Range("K26").Value = 6.472491909385
$ws.Range("L26").Value = -2.662721893491
$ws.Range("M26").Value = -24.885844748858
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 48
$ws.Range("K28").Value = 41.176470588235
$ws.Range("L28").Value = 29.729729729729
$ws.Range("N29").Value = -69.230769230769
$ws.Range("N30").Value = -69.230769230769

# --- Convert C18 from text placeholder back to a number; fix its style to match the numeric column style ---
$ws.Range("F18").Copy()
$ws.Range("C18").PasteSpecial(-4122)

# --- Cells that become "no data" text placeholders (shared strings "0" / "***.*") ---
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "***.*"
$ws.Range("G22").Value = "'0"
$ws.Range("H22").Value = "***.*"
$ws.Range("G23").Value = "'0"
$ws.Range("H23").Value = "***.*"
$ws.Range("C27").Value = "'0"
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "***.*"
$ws.Range("D28").Value = "'0"
$ws.Range("E28").Value = "***.*"

# --- Restore style 13 (placeholder style) on the cells just converted to text; copy formats from row 14 reference cells ---
$ws.Range("C14:E14").Copy()
$ws.Range("C15:E15").PasteSpecial(-4122)
$ws.Range("C27:E27").PasteSpecial(-4122)
$ws.Range("C14:D14").Copy()
$ws.Range("G22:H22").PasteSpecial(-4122)
$ws.Range("G23:H23").PasteSpecial(-4122)
$ws.Range("D14:E14").Copy()
$ws.Range("D28:E28").PasteSpecial(-4122)
